# Gallery_PanelNode_Repeaters_ProPanels.xlsx
# "Added Test data for Spain Zettler Market"
#
# The workbook currently ends with the "Italy" sheet (7th / last tab,
# currently the active tab). We duplicate it to create a new "Spain"
# sheet immediately after it, fill in the market name / part number,
# and move the active-tab/selection state from Italy onto the new sheet
# (mirroring what Excel does when a user right-clicks a tab -> Move or
# Copy... -> Create a copy, renames it, and edits the two data cells).

$wb = $excel.ActiveWorkbook

# --- duplicate the last sheet ("Italy") to the end of the tab strip ---
$sheetCount = $wb.Worksheets.Count
$sourceSheet = $wb.Worksheets.Item($sheetCount)
$sourceSheet.Copy([System.Reflection.Missing]::Value, $sourceSheet)

$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# --- fill in the market-specific values (same layout as every other market tab) ---
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2034"

# the row that holds the part number code wraps onto two lines for Spain,
# same as it visually does for the other tabs once the value is entered
$spain.Rows.Item(3).RowHeight = 28.8
$spain.Rows.Item(4).RowHeight = 28.8
$spain.Rows.Item(5).RowHeight = 28.8

# columns sized to fit the (slightly different) Spain content
$spain.Columns.Item(1).ColumnWidth = 24.276041666666668
$spain.Columns.Item(2).ColumnWidth = 15.944010416666666
$spain.Columns.Item(4).ColumnWidth = 17.385416666666668

# --- move the "active tab" selection state from Italy onto the new Spain tab ---
$italy = $wb.Worksheets.Item($sheetCount)
$italy.Range("A1:D21").Select() | Out-Null

$spain.Activate() | Out-Null
$spain.Range("F14").Select() | Out-Null
